# CNCBBS.xlsx update: insert a new "45412 / 45443 / 45473" data block after the
# existing 37-row block (old rows 39-75 shift down to 42-78), and append the
# same 3-row block at the very end (new rows 79-81).
#
# Net effect: sheet grows from A1:F75 to A1:F81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for 3 new rows right after the first data block (row 38) ---
# This pushes the former rows 39:75 down to 42:78 and bumps the dimension/used range.
$ws.Rows("39:41").Insert()

# The rows that just landed at 42:44 are a verbatim copy of the old 39:41, so their
# formatting (style "2": bordered, centered, custom date numfmt) is exactly what the
# brand-new rows should look like too. Clone just that formatting (columns A:F only)
# into the freshly inserted, still-blank rows.
$ws.Range("A42:F44").Copy()
$ws.Range("A39:F41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: populate the 3 newly inserted rows (39:41) ---
$ws.Range("A39").Value = 45412
$ws.Range("B39:E39").Value = 42842522000000
$ws.Range("F39").Value = 0

$ws.Range("A40").Value = 45443
$ws.Range("B40:E40").Value = 42953758000000
$ws.Range("F40").Value = 0

$ws.Range("A41").Value = 45473
$ws.Range("B41:E41").Value = 43654679000000
$ws.Range("F41").Value = 0

# --- Step 3: append the same 3-row block at the end of the sheet (new rows 79:81) ---
# Clone formatting from the last existing data row (now row 38, since it's outside the
# shifted block) down into the new trailing rows.
$ws.Range("A38:F38").Copy()
$ws.Range("A79:F81").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A79").Value = 45412
$ws.Range("B79:E79").Value = 42842522000000
$ws.Range("F79").Value = 0

$ws.Range("A80").Value = 45443
$ws.Range("B80:E80").Value = 42953758000000
$ws.Range("F80").Value = 0

$ws.Range("A81").Value = 45473
$ws.Range("B81:E81").Value = 43654679000000
$ws.Range("F81").Value = 0
